$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 100623597
$ws.Range("B2").Value = "Yusuf Pramudya"
$ws.Range("C2").Value = "Bank Mandiri"
$ws.Range("D2").Value = "Yusuf Pramudya Hutama"
$ws.Range("E2").Value = "'1560022694287"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 2692500
$ws.Range("H2").Value = 2692500
$ws.Range("I2").Value = 45658
$ws.Range("J2").Value = 45672
$ws.Range("K2").Value = "Yusuf_Pramudya"

# Row 3
$ws.Range("A3").Value = 100623598
$ws.Range("B3").Value = "Fajar Dafa"
$ws.Range("C3").Value = "Bank BCA"
$ws.Range("D3").Value = "Fajar Dafa' Supriyanto "
$ws.Range("E3").Value = "'5776390247"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = 102500
$ws.Range("G3").Value = 2080000
$ws.Range("H3").Value = 1977500
$ws.Range("I3").Value = 45658
$ws.Range("J3").Value = 45672
$ws.Range("K3").Value = "Fajar_Dafa"
